# Quarterly indexing esoteric bug-fix operation
# Shifts the date serials in column A (rows 2-73) from the 1st of the
# quarter-start month to the 15th of the following month (mid-quarter),
# i.e. each original serial is advanced by ~44-45 days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @{
    2  = 25614
    3  = 25614
    4  = 25614
    5  = 25614
    6  = 25614
    7  = 25614
    8  = 25614
    9  = 25614
    10 = 25614
    11 = 25614
    12 = 39401
    13 = 39493
    14 = 39583
    15 = 39675
    16 = 39767
    17 = 39859
    18 = 39948
    19 = 40040
    20 = 40132
    21 = 40224
    22 = 40313
    23 = 40405
    24 = 40497
    25 = 40589
    26 = 40678
    27 = 40770
    28 = 40862
    29 = 40954
    30 = 41044
    31 = 41136
    32 = 41228
    33 = 41320
    34 = 41409
    35 = 41501
    36 = 41593
    37 = 41685
    38 = 41774
    39 = 41866
    40 = 41958
    41 = 42050
    42 = 42139
    43 = 42231
    44 = 42323
    45 = 42415
    46 = 42505
    47 = 42597
    48 = 42689
    49 = 42781
    50 = 42870
    51 = 42962
    52 = 43054
    53 = 43146
    54 = 43235
    55 = 43327
    56 = 43419
    57 = 43511
    58 = 43600
    59 = 43692
    60 = 43784
    61 = 43876
    62 = 43966
    63 = 44058
    64 = 44150
    65 = 44242
    66 = 44331
    67 = 44423
    68 = 44515
    69 = 44607
    70 = 44696
    71 = 44788
    72 = 44880
    73 = 44972
}

foreach ($row in $newDates.Keys) {
    $ws.Cells.Item($row, 1).Value = $newDates[$row]
}
